$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date number format used by the existing "Date Delivered to Xin" column (F)
$dateFmt = $ws.Range("F2").NumberFormat

# --- Append the newly launched sensors (rows 22-27) ---
$newRows = @(
    @{ Row = 22; A = "PT9-M8U";  B = "645156ce47e3b51c2399761a"; C = "645156cf1ca5eb000db782ae"; D = 44.205145999999999; E = -107.92413500000001; F = 45470 }
    @{ Row = 23; A = "WS7-7R8";  B = "65c14dc29deec5000eb8c4ec"; C = "65c14dc29deec5000c8f3edc"; D = 40.542369999999998;  E = -102.2972;            F = 45470 }
    @{ Row = 24; A = "WS47-6P9"; B = "6644d5387d2393000b8b9e7f"; C = "6644d5387d2393000b8b9e80"; D = 41.062220000000003;  E = -101.95968999999999; F = 45470 }
    @{ Row = 25; A = "PT16-VMF"; B = "6446abe57e1943000ed84f1a"; C = "6446abe5d4c889000e636142"; D = 40.774619999999999;  E = -103.03641;           F = 45470 }
    @{ Row = 26; A = "WS5-G2N";  B = "65c13e5463c11f000d32382c"; C = "65c13e559deec5000eb8c4ea"; D = 41.105719999999998;  E = -101.69614;           F = 45470 }
    @{ Row = 27; A = "WS3-NYM";  B = "65c13224b78889000c09cfb4"; C = "65c1322429e691000e783bc4"; D = 40.581780000000002;  E = -102.03254;           F = 45470 }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("F$r").NumberFormat = $dateFmt
}

# --- Fill the "delivered by" column (G), newest batch to oldest ---
$ws.Range("G22").Value = "Mark Bjornstead"
foreach ($r in 23..27) { $ws.Range("G$r").Value = "Kreg Vollmer" }
foreach ($r in @(20, 19, 18)) { $ws.Range("G$r").Value = "Leah Wimmer" }
foreach ($r in @(16, 15, 14, 13)) { $ws.Range("G$r").Value = "Stephen Rose" }
foreach ($r in @(11, 10, 9, 8)) { $ws.Range("G$r").Value = "Alyssa Brewer" }
foreach ($r in @(6, 5, 4, 3, 2)) { $ws.Range("G$r").Value = "Andy Stieger" }

$ws.Range("E32").Select()
